# Insert a new weekly price record as row 136 on Sheet1.
# This shifts the existing rows 136-159 down to 137-160 (matching the
# new dimension A1:R160) and fills the newly inserted row with the
# latest weekly observation for "Espinaca" at Terminal La Palmera de La
# Serena.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 136, shifting rows 136:159 down
# to 137:160 and carrying the existing row formatting (date style, etc.)
$ws.Rows.Item(136).Insert()

$ws.Cells.Item(136, 1).Value = 8
$ws.Cells.Item(136, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(136, 3).Value = "Coquimbo"
$ws.Cells.Item(136, 4).Value = 44476
$ws.Cells.Item(136, 5).Value = 4
$ws.Cells.Item(136, 6).Value = 100112012
$ws.Cells.Item(136, 7).Value = "Espinaca"
$ws.Cells.Item(136, 8).Value = "Sin especificar"
$ws.Cells.Item(136, 9).Value = "Primera"
$ws.Cells.Item(136, 10).Value = 2800
$ws.Cells.Item(136, 11).Value = 400
$ws.Cells.Item(136, 12).Value = 500
$ws.Cells.Item(136, 13).Value = 450
$ws.Cells.Item(136, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(136, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(136, 16).Value = 900
$ws.Cells.Item(136, 17).Value = 0.5
$ws.Cells.Item(136, 18).Value = "Hortaliza"

Write-Output "Inserted row 136; dimension now $($ws.UsedRange.Rows.Count) rows"
